$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (serial date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$newRows = @(
    @(44441, 22, 176, 93.11528836640865),
    @(44442, 47, 187, 98.9349938893092),
    @(44443, 24, 183, 96.81873733552719),
    @(44444, 26, 186, 98.40592975086369),
    @(44445, 55, 215, 113.7487897657833),
    @(44446, 17, 205, 108.4581483813283),
    @(44447, 9, 200, 105.8128276891007),
    @(44448, 18, 196, 103.6965711353187)
)

$startRow = 367
$templateCellA = $ws.Cells.Item($startRow - 1, 1)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $templateCellA.Copy($cellA)
    $cellA.Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
